# Applies the data_selected.xlsx falling-data edit:
#  - 3 new accelerometer/gyroscope sample rows are inserted at the top of the
#    data block (rows 2-4), the remaining original samples shift down,
#    and the sheet keeps exactly one extra row overall (dimension -> A1:H21).
#  - timestamp/label columns (A,B) stay tied to their row position (0,100,200,...)
#    and a new trailing sample row (timestamp 1900) is appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "falling"
$ws.Cells.Item(2, 3).Value = -2.669419974088668
$ws.Cells.Item(2, 4).Value = 9.347340643405914
$ws.Cells.Item(2, 5).Value = -0.05590170621871929
$ws.Cells.Item(2, 6).Value = 0.0371100641787052
$ws.Cells.Item(2, 7).Value = -0.0024434609804302
$ws.Cells.Item(2, 8).Value = 0.0224492978304624

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "falling"
$ws.Cells.Item(3, 3).Value = -2.789929866790771
$ws.Cells.Item(3, 4).Value = 9.389312267303467
$ws.Cells.Item(3, 5).Value = -0.0143058076500895
$ws.Cells.Item(3, 6).Value = 0.0600175112485885
$ws.Cells.Item(3, 7).Value = -0.0474947728216648
$ws.Cells.Item(3, 8).Value = 0.0047342055477201

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "falling"
$ws.Cells.Item(4, 3).Value = -3.052737355232238
$ws.Cells.Item(4, 4).Value = 9.12702190876007
$ws.Cells.Item(4, 5).Value = -0.5663906224071975
$ws.Cells.Item(4, 6).Value = 0.0438295826315879
$ws.Cells.Item(4, 7).Value = -0.0827722400426864
$ws.Cells.Item(4, 8).Value = 0.0759000033140182

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "falling"
$ws.Cells.Item(5, 3).Value = -3.195896685123444
$ws.Cells.Item(5, 4).Value = 8.818168640136719
$ws.Cells.Item(5, 5).Value = -1.334951654076576
$ws.Cells.Item(5, 6).Value = -0.0627664029598236
$ws.Cells.Item(5, 7).Value = -0.5068654417991638
$ws.Cells.Item(5, 8).Value = -0.0064140851609408

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "falling"
$ws.Cells.Item(6, 3).Value = -4.037624061107635
$ws.Cells.Item(6, 4).Value = 8.026161462068558
$ws.Cells.Item(6, 5).Value = -2.41215243935585
$ws.Cells.Item(6, 6).Value = -0.0687223374843597
$ws.Cells.Item(6, 7).Value = -0.2246456891298294
$ws.Cells.Item(6, 8).Value = 0.1652390509843826

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "falling"
$ws.Cells.Item(7, 3).Value = -3.571272850036621
$ws.Cells.Item(7, 4).Value = 6.589316844940186
$ws.Cells.Item(7, 5).Value = -3.011023998260498
$ws.Cells.Item(7, 6).Value = -0.2229658216238021
$ws.Cells.Item(7, 7).Value = -0.113315500319004
$ws.Cells.Item(7, 8).Value = 0.2180788964033126

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "falling"
$ws.Cells.Item(8, 3).Value = -5.014195352792741
$ws.Cells.Item(8, 4).Value = 5.515790849924088
$ws.Cells.Item(8, 5).Value = -4.485763758420944
$ws.Cells.Item(8, 6).Value = 0.3182607889175415
$ws.Cells.Item(8, 7).Value = 0.6488916277885437
$ws.Cells.Item(8, 8).Value = 0.4970915913581848

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "falling"
$ws.Cells.Item(9, 3).Value = -6.389075994491578
$ws.Cells.Item(9, 4).Value = 5.315201640129089
$ws.Cells.Item(9, 5).Value = -2.994861543178557
$ws.Cells.Item(9, 6).Value = 0.3697261810302734
$ws.Cells.Item(9, 7).Value = 1.551903128623962
$ws.Cells.Item(9, 8).Value = -0.230448916554451

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "falling"
$ws.Cells.Item(10, 3).Value = -4.530004173517225
$ws.Cells.Item(10, 4).Value = 6.061075717210771
$ws.Cells.Item(10, 5).Value = -1.65171818435192
$ws.Cells.Item(10, 6).Value = 0.451123982667923
$ws.Cells.Item(10, 7).Value = 0.8857545852661133
$ws.Cells.Item(10, 8).Value = -1.012661814689636

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "falling"
$ws.Cells.Item(11, 3).Value = -4.840012192726136
$ws.Cells.Item(11, 4).Value = 7.050750851631165
$ws.Cells.Item(11, 5).Value = 2.538701653480534
$ws.Cells.Item(11, 6).Value = -0.1492038369178772
$ws.Cells.Item(11, 7).Value = -1.078635334968567
$ws.Cells.Item(11, 8).Value = -4.738939762115479

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "falling"
$ws.Cells.Item(12, 3).Value = 39.4459085166455
$ws.Cells.Item(12, 4).Value = 14.87950980663301
$ws.Cells.Item(12, 5).Value = 20.98658950626853
$ws.Cells.Item(12, 6).Value = -1.247386813163757
$ws.Cells.Item(12, 7).Value = -1.642922043800354
$ws.Cells.Item(12, 8).Value = 0.3756821155548095

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "falling"
$ws.Cells.Item(13, 3).Value = 20.97534060478203
$ws.Cells.Item(13, 4).Value = 11.63514316082
$ws.Cells.Item(13, 5).Value = 12.98352101445195
$ws.Cells.Item(13, 6).Value = 2.319607973098755
$ws.Cells.Item(13, 7).Value = -0.2393064647912979
$ws.Cells.Item(13, 8).Value = -2.242333650588989

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "falling"
$ws.Cells.Item(14, 3).Value = 1.842506676912319
$ws.Cells.Item(14, 4).Value = 6.827371656894682
$ws.Cells.Item(14, 5).Value = 4.018280878663065
$ws.Cells.Item(14, 6).Value = 0.1510364264249801
$ws.Cells.Item(14, 7).Value = 2.41413950920105
$ws.Cells.Item(14, 8).Value = 0.4211915731430053

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "falling"
$ws.Cells.Item(15, 3).Value = 5.506411552429201
$ws.Cells.Item(15, 4).Value = 6.245316505432129
$ws.Cells.Item(15, 5).Value = 3.238075017929073
$ws.Cells.Item(15, 6).Value = 0.2503020465373993
$ws.Cells.Item(15, 7).Value = -1.337489485740662
$ws.Cells.Item(15, 8).Value = 0.0983493030071258

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "falling"
$ws.Cells.Item(16, 3).Value = 6.049000829458237
$ws.Cells.Item(16, 4).Value = 7.41145807504654
$ws.Cells.Item(16, 5).Value = 3.085425972938538
$ws.Cells.Item(16, 6).Value = -0.4132503271102905
$ws.Cells.Item(16, 7).Value = 0.2756529450416565
$ws.Cells.Item(16, 8).Value = -0.5613851547241211

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "falling"
$ws.Cells.Item(17, 3).Value = 5.151385545730594
$ws.Cells.Item(17, 4).Value = 9.354082107543944
$ws.Cells.Item(17, 5).Value = 3.488138377666476
$ws.Cells.Item(17, 6).Value = -0.1337794959545135
$ws.Cells.Item(17, 7).Value = 0.3077233731746673
$ws.Cells.Item(17, 8).Value = -0.0181732401251792

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "falling"
$ws.Cells.Item(18, 3).Value = 3.843811631202696
$ws.Cells.Item(18, 4).Value = 9.246257454156877
$ws.Cells.Item(18, 5).Value = 2.557389497756956
$ws.Cells.Item(18, 6).Value = 0.0068722339347004
$ws.Cells.Item(18, 7).Value = -0.2590068578720093
$ws.Cells.Item(18, 8).Value = 0.2180788964033126

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "falling"
$ws.Cells.Item(19, 3).Value = 4.352240324020387
$ws.Cells.Item(19, 4).Value = 7.618153929710387
$ws.Cells.Item(19, 5).Value = 3.785528540611267
$ws.Cells.Item(19, 6).Value = -0.0940732508897781
$ws.Cells.Item(19, 7).Value = -0.0920879393815994
$ws.Cells.Item(19, 8).Value = 0.011148290708661

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "falling"
$ws.Cells.Item(20, 3).Value = 3.598823010921479
$ws.Cells.Item(20, 4).Value = 8.679397195577621
$ws.Cells.Item(20, 5).Value = 3.499203532934189
$ws.Cells.Item(20, 6).Value = -0.0429132841527462
$ws.Cells.Item(20, 7).Value = -0.0099265603348612
$ws.Cells.Item(20, 8).Value = -0.2417499274015426

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "falling"
$ws.Cells.Item(21, 3).Value = 3.266197681427002
$ws.Cells.Item(21, 4).Value = 7.760588467121124
$ws.Cells.Item(21, 5).Value = 3.241497814655304
$ws.Cells.Item(21, 6).Value = 0.0386372283101081
$ws.Cells.Item(21, 7).Value = 0.0759000033140182
$ws.Cells.Item(21, 8).Value = 0.0108428578823804

